# Auto-generated: updates FFXIV Ifrit market-price derived columns (H:N)
# for specific leve rows across multiple sheets, per scheduled market refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 299.75
$ws.Range("I2").Value = 333
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 333
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = -220
$ws.Range("N2").Value = -426

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 197.2
$ws.Range("I38").Value = 69
$ws.Range("J38").Value = 549.75
$ws.Range("K38").Value = 207
$ws.Range("L38").Value = 1649.25
$ws.Range("M38").Value = 165
$ws.Range("N38").Value = -2393.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 457.875
$ws.Range("I39").Value = 84.59999999999999
$ws.Range("J39").Value = 1080
$ws.Range("K39").Value = 253.8
$ws.Range("L39").Value = 3240
$ws.Range("M39").Value = 42.20000000000002
$ws.Range("N39").Value = -3832

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1056
$ws.Range("I40").Value = 1001.53845
$ws.Range("J40").Value = 1157.1428
$ws.Range("K40").Value = 1001.53845
$ws.Range("L40").Value = 1157.1428
$ws.Range("M40").Value = -826.53845
$ws.Range("N40").Value = -1507.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 158.6
$ws.Range("I61").Value = 158.6
$ws.Range("K61").Value = 475.8
$ws.Range("M61").Value = -303.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 710.25
$ws.Range("I129").Value = 409.4
$ws.Range("J129").Value = 925.1429000000001
$ws.Range("K129").Value = 1228.2
$ws.Range("L129").Value = 2775.4287
$ws.Range("M129").Value = 3771.8
$ws.Range("N129").Value = -12775.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4331.577
$ws.Range("I132").Value = 4443.5415
$ws.Range("J132").Value = 2988
$ws.Range("K132").Value = 13330.6245
$ws.Range("L132").Value = 8964
$ws.Range("M132").Value = -10800.6245
$ws.Range("N132").Value = -14024

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 50002000
$ws.Range("I137").Value = 1749.3334
$ws.Range("J137").Value = 125002376
$ws.Range("K137").Value = 5248.0002
$ws.Range("L137").Value = 375007128
$ws.Range("M137").Value = -2698.0002
$ws.Range("N137").Value = -375012228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 910.6667
$ws.Range("I25").Value = 910.6667
$ws.Range("K25").Value = 910.6667
$ws.Range("M25").Value = -508.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 4366.6665
$ws.Range("I28").Value = 1280
$ws.Range("J28").Value = 19800
$ws.Range("K28").Value = 1280
$ws.Range("L28").Value = 19800
$ws.Range("M28").Value = -1088
$ws.Range("N28").Value = -20184

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7409969.5
$ws.Range("I61").Value = 8549400
$ws.Range("J61").Value = 3675
$ws.Range("K61").Value = 8549400
$ws.Range("L61").Value = 3675
$ws.Range("M61").Value = -8549188
$ws.Range("N61").Value = -4099

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 511.75674
$ws.Range("I97").Value = 490.92593
$ws.Range("K97").Value = 490.92593
$ws.Range("M97").Value = 5.074070000000006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 4366.6665
$ws.Range("I99").Value = 1280
$ws.Range("J99").Value = 19800
$ws.Range("K99").Value = 1280
$ws.Range("L99").Value = 19800
$ws.Range("M99").Value = 1715
$ws.Range("N99").Value = -25790

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 28978
$ws.Range("J113").Value = 28978
$ws.Range("L113").Value = 28978
$ws.Range("N113").Value = -37656

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7409969.5
$ws.Range("I136").Value = 8549400
$ws.Range("J136").Value = 3675
$ws.Range("K136").Value = 25648200
$ws.Range("L136").Value = 11025
$ws.Range("M136").Value = -25645650
$ws.Range("N136").Value = -16125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1730.0435
$ws.Range("I31").Value = 928.8823
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 928.8823
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -633.8823
$ws.Range("N31").Value = -4590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1730.0435
$ws.Range("I34").Value = 928.8823
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 928.8823
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -726.8823
$ws.Range("N34").Value = -4404

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4390.8335
$ws.Range("I58").Value = 1843.7142
$ws.Range("J58").Value = 5005.6553
$ws.Range("K58").Value = 1843.7142
$ws.Range("L58").Value = 5005.6553
$ws.Range("M58").Value = -1640.7142
$ws.Range("N58").Value = -5411.6553

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 31000
$ws.Range("I63").Value = 31000
$ws.Range("K63").Value = 31000
$ws.Range("M63").Value = -30314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 31000
$ws.Range("I66").Value = 31000
$ws.Range("K66").Value = 93000
$ws.Range("M66").Value = -89568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4390.8335
$ws.Range("I136").Value = 1843.7142
$ws.Range("J136").Value = 5005.6553
$ws.Range("K136").Value = 5531.142599999999
$ws.Range("L136").Value = 15016.9659
$ws.Range("M136").Value = -2981.142599999999
$ws.Range("N136").Value = -20116.9659

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 23809952
$ws.Range("I5").Value = 408.11765
$ws.Range("J5").Value = 125000510
$ws.Range("K5").Value = 1224.35295
$ws.Range("L5").Value = 375001530
$ws.Range("M5").Value = -1112.35295
$ws.Range("N5").Value = -375001754

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2750.8223
$ws.Range("I113").Value = 563.6875
$ws.Range("J113").Value = 3957.5173
$ws.Range("K113").Value = 1691.0625
$ws.Range("L113").Value = 11872.5519
$ws.Range("M113").Value = 478.9375
$ws.Range("N113").Value = -16212.5519

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7148.853
$ws.Range("I131").Value = 9051.429
$ws.Range("J131").Value = 6655.593
$ws.Range("K131").Value = 27154.287
$ws.Range("L131").Value = 19966.779
$ws.Range("M131").Value = -22114.287
$ws.Range("N131").Value = -30046.779

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 23809952
$ws.Range("I135").Value = 408.11765
$ws.Range("J135").Value = 125000510
$ws.Range("K135").Value = 3673.05885
$ws.Range("L135").Value = 1125004590
$ws.Range("M135").Value = -1138.05885
$ws.Range("N135").Value = -1125009660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 7350
$ws.Range("I141").Value = 7350
$ws.Range("K141").Value = 22050
$ws.Range("M141").Value = -16870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 23892.25
$ws.Range("J32").Value = 23892.25
$ws.Range("L32").Value = 23892.25
$ws.Range("N32").Value = -24484.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2132.9443
$ws.Range("I80").Value = 2319.2307
$ws.Range("J80").Value = 1648.6
$ws.Range("K80").Value = 2319.2307
$ws.Range("L80").Value = 1648.6
$ws.Range("M80").Value = -1321.2307
$ws.Range("N80").Value = -3644.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2132.9443
$ws.Range("I83").Value = 2319.2307
$ws.Range("J83").Value = 1648.6
$ws.Range("K83").Value = 11596.1535
$ws.Range("L83").Value = 8243
$ws.Range("M83").Value = -6604.1535
$ws.Range("N83").Value = -18227

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1413
$ws.Range("I82").Value = 1365.8823
$ws.Range("J82").Value = 1546.5
$ws.Range("K82").Value = 1365.8823
$ws.Range("L82").Value = 1546.5
$ws.Range("M82").Value = -1004.8823
$ws.Range("N82").Value = -2268.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1413
$ws.Range("I85").Value = 1365.8823
$ws.Range("J85").Value = 1546.5
$ws.Range("K85").Value = 1365.8823
$ws.Range("L85").Value = 1546.5
$ws.Range("M85").Value = -117.8823
$ws.Range("N85").Value = -4042.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 13370
$ws.Range("J103").Value = 13370
$ws.Range("L103").Value = 13370
$ws.Range("N103").Value = -15714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 16975
$ws.Range("J74").Value = 19466.666
$ws.Range("L74").Value = 19466.666
$ws.Range("N74").Value = -21338.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 16975
$ws.Range("J77").Value = 19466.666
$ws.Range("L77").Value = 58399.99800000001
$ws.Range("N77").Value = -67759.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4907.6665
$ws.Range("I81").Value = 2237.75
$ws.Range("J81").Value = 6242.625
$ws.Range("K81").Value = 4475.5
$ws.Range("L81").Value = 12485.25
$ws.Range("M81").Value = -3414.5
$ws.Range("N81").Value = -14607.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4907.6665
$ws.Range("I84").Value = 2237.75
$ws.Range("J84").Value = 6242.625
$ws.Range("K84").Value = 22377.5
$ws.Range("L84").Value = 62426.25
$ws.Range("M84").Value = -17073.5
$ws.Range("N84").Value = -73034.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 46000
$ws.Range("J93").Value = 46000
$ws.Range("L93").Value = 46000
$ws.Range("N93").Value = -50992
